$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 53.333332
$ws.Range("I9").Value = 30
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 30
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 139
$ws.Range("N9").Value = -438

$ws.Range("H98").Value = 1020.8889
$ws.Range("I98").Value = 773.05
$ws.Range("K98").Value = 773.05
$ws.Range("M98").Value = 724.95

$ws.Range("H116").Value = 4423.375
$ws.Range("I116").Value = 3736.75
$ws.Range("J116").Value = 5110
$ws.Range("K116").Value = 3736.75
$ws.Range("L116").Value = 5110
$ws.Range("M116").Value = -294.75
$ws.Range("N116").Value = -11994

$ws.Range("H122").Value = 1020.8889
$ws.Range("I122").Value = 773.05
$ws.Range("K122").Value = 2319.15
$ws.Range("M122").Value = 130.8500000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8234.754999999999
$ws.Range("I32").Value = 7520.0596
$ws.Range("J32").Value = 12522.929
$ws.Range("K32").Value = 7520.0596
$ws.Range("L32").Value = 12522.929
$ws.Range("M32").Value = -7233.0596
$ws.Range("N32").Value = -13096.929

$ws.Range("H45").Value = 9976
$ws.Range("I45").Value = 11125.632
$ws.Range("J45").Value = 7990.273
$ws.Range("K45").Value = 11125.632
$ws.Range("L45").Value = 7990.273
$ws.Range("M45").Value = -10748.632
$ws.Range("N45").Value = -8744.273000000001

$ws.Range("H97").Value = 28601598
$ws.Range("I97").Value = 43479668
$ws.Range("J97").Value = 85298.25
$ws.Range("K97").Value = 43479668
$ws.Range("L97").Value = 85298.25
$ws.Range("M97").Value = -43479172
$ws.Range("N97").Value = -86290.25

$ws.Range("H105").Value = 39974
$ws.Range("J105").Value = 39974
$ws.Range("L105").Value = 39974
$ws.Range("N105").Value = -46962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4124.033
$ws.Range("I20").Value = 4344.05
$ws.Range("J20").Value = 3684
$ws.Range("K20").Value = 4344.05
$ws.Range("L20").Value = 3684
$ws.Range("M20").Value = -4097.05
$ws.Range("N20").Value = -4178

$ws.Range("H86").Value = 2205.9048
$ws.Range("I86").Value = 1908.75
$ws.Range("J86").Value = 2602.111
$ws.Range("K86").Value = 1908.75
$ws.Range("L86").Value = 2602.111
$ws.Range("M86").Value = -785.75
$ws.Range("N86").Value = -4848.111

$ws.Range("H89").Value = 2205.9048
$ws.Range("I89").Value = 1908.75
$ws.Range("J89").Value = 2602.111
$ws.Range("K89").Value = 9543.75
$ws.Range("L89").Value = 13010.555
$ws.Range("M89").Value = -3927.75
$ws.Range("N89").Value = -24242.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2221
$ws.Range("I16").Value = 526.25
$ws.Range("J16").Value = 9000
$ws.Range("K16").Value = 526.25
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -239.25
$ws.Range("N16").Value = -9574

$ws.Range("H31").Value = 4138.1
$ws.Range("I31").Value = 2689.3684
$ws.Range("J31").Value = 6640.4546
$ws.Range("K31").Value = 2689.3684
$ws.Range("L31").Value = 6640.4546
$ws.Range("M31").Value = -2394.3684
$ws.Range("N31").Value = -7230.4546

$ws.Range("H34").Value = 4138.1
$ws.Range("I34").Value = 2689.3684
$ws.Range("J34").Value = 6640.4546
$ws.Range("K34").Value = 2689.3684
$ws.Range("L34").Value = 6640.4546
$ws.Range("M34").Value = -2487.3684
$ws.Range("N34").Value = -7044.4546

$ws.Range("H92").Value = 36000
$ws.Range("J92").Value = 36000
$ws.Range("L92").Value = 36000
$ws.Range("N92").Value = -40992

$ws.Range("H99").Value = 9234.306
$ws.Range("I99").Value = 6499.6924
$ws.Range("J99").Value = 12325.608
$ws.Range("K99").Value = 6499.6924
$ws.Range("L99").Value = 12325.608
$ws.Range("M99").Value = -5001.6924
$ws.Range("N99").Value = -15321.608

$ws.Range("H105").Value = 1665.3846
$ws.Range("I105").Value = 1415.1
$ws.Range("J105").Value = 2499.6667
$ws.Range("K105").Value = 1415.1
$ws.Range("L105").Value = 2499.6667
$ws.Range("M105").Value = 331.9000000000001
$ws.Range("N105").Value = -5993.6667

$ws.Range("H113").Value = 2221
$ws.Range("I113").Value = 526.25
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 526.25
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = 1643.75
$ws.Range("N113").Value = -13340

$ws.Range("H126").Value = 9234.306
$ws.Range("I126").Value = 6499.6924
$ws.Range("J126").Value = 12325.608
$ws.Range("K126").Value = 19499.0772
$ws.Range("L126").Value = 36976.824
$ws.Range("M126").Value = -17029.0772
$ws.Range("N126").Value = -41916.824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49.333332
$ws.Range("J12").Value = 39.11111
$ws.Range("L12").Value = 117.33333
$ws.Range("N12").Value = -463.33333

$ws.Range("H131").Value = 3268963
$ws.Range("I131").Value = 4202523.5
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 12607570.5
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = -12602530.5
$ws.Range("N131").Value = -14580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5699
$ws.Range("I70").Value = 5998
$ws.Range("J70").Value = 5400
$ws.Range("K70").Value = 5998
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -5728
$ws.Range("N70").Value = -5940

$ws.Range("H73").Value = 5699
$ws.Range("I73").Value = 5998
$ws.Range("J73").Value = 5400
$ws.Range("K73").Value = 5998
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -5062
$ws.Range("N73").Value = -7272

$ws.Range("H95").Value = 101655.14
$ws.Range("J95").Value = 101655.14
$ws.Range("L95").Value = 101655.14
$ws.Range("N95").Value = -107147.14

$ws.Range("H102").Value = 3681.5696
$ws.Range("I102").Value = 3006.1453
$ws.Range("J102").Value = 6144.8823
$ws.Range("K102").Value = 3006.1453
$ws.Range("L102").Value = 6144.8823
$ws.Range("M102").Value = -1384.1453
$ws.Range("N102").Value = -9388.882300000001

$ws.Range("H113").Value = 5999
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5999
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10339

$ws.Range("H126").Value = 4027.4119
$ws.Range("I126").Value = 3782.1025
$ws.Range("J126").Value = 4824.6665
$ws.Range("K126").Value = 11346.3075
$ws.Range("L126").Value = 14473.9995
$ws.Range("M126").Value = -8876.307499999999
$ws.Range("N126").Value = -19413.9995

$ws.Range("H136").Value = 69977.25
$ws.Range("J136").Value = 69977.25
$ws.Range("L136").Value = 209931.75
$ws.Range("N136").Value = -215031.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9001
$ws.Range("I2").Value = 9000
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = -8888

$ws.Range("H40").Value = 2773.9375
$ws.Range("I40").Value = 2773
$ws.Range("J40").Value = 2776.75
$ws.Range("K40").Value = 2773
$ws.Range("L40").Value = 2776.75
$ws.Range("M40").Value = -2637
$ws.Range("N40").Value = -3048.75

$ws.Range("H45").Value = 14749.5
$ws.Range("J45").Value = 8000
$ws.Range("L45").Value = 8000
$ws.Range("N45").Value = -8814

$ws.Range("H46").Value = 2719.6316
$ws.Range("I46").Value = 1022
$ws.Range("J46").Value = 3709.9167
$ws.Range("K46").Value = 1022
$ws.Range("L46").Value = 3709.9167
$ws.Range("M46").Value = -834
$ws.Range("N46").Value = -4085.9167

$ws.Range("H61").Value = 3419.0952
$ws.Range("I61").Value = 3419.0952
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3419.0952
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3217.0952
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 3419.0952
$ws.Range("I113").Value = 3419.0952
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3419.0952
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1249.0952
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 54349.668
$ws.Range("J63").Value = 54349.668
$ws.Range("L63").Value = 54349.668
$ws.Range("N63").Value = -55597.668

$ws.Range("H66").Value = 54349.668
$ws.Range("J66").Value = 54349.668
$ws.Range("L66").Value = 163049.004
$ws.Range("N66").Value = -169289.004

$ws.Range("H107").Value = 1227.3462
$ws.Range("I107").Value = 1204.381
$ws.Range("J107").Value = 1323.8
$ws.Range("K107").Value = 3613.143
$ws.Range("L107").Value = 3971.4
$ws.Range("M107").Value = -1693.143
$ws.Range("N107").Value = -7811.4
